# Natmi following Dr Hou advice
# Update C3-Cd46 LR-pair table: add "ECs" as a Sending cluster group
# (rows 2-4), shift the FAPs/sCs groups down, and append a new sCs
# sending-cluster group (rows 8-10) with refreshed edge-expression stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'ECs'
$ws.Range('B2').Value = 'C3'
$ws.Range('C2').Value = 'Cd46'
$ws.Range('D2').Value = 'ECs'
$ws.Range('E2').Value = 2
$ws.Range('F2').Value = 0.6666666666666666
$ws.Range('G2').Value = 34.739995
$ws.Range('H2').Value = 104.219985
$ws.Range('I2').Value = 0.1827267341390226
$ws.Range('J2').Value = 0.1827267341390226
$ws.Range('K2').Value = 3
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = 10.876357
$ws.Range('N2').Value = 32.629071
$ws.Range('O2').Value = 0.6153167079818751
$ws.Range('P2').Value = 0.615316707981875
$ws.Range('Q2').Value = 377.8445877982149
$ws.Range('R2').Value = 3400.601290183934
$ws.Range('S2').Value = 0.1124348125107027
$ws.Range('T2').Value = 0.1124348125107027

# Row 3
$ws.Range('A3').Value = 'ECs'
$ws.Range('B3').Value = 'C3'
$ws.Range('C3').Value = 'Cd46'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('E3').Value = 2
$ws.Range('F3').Value = 0.6666666666666666
$ws.Range('G3').Value = 34.739995
$ws.Range('H3').Value = 104.219985
$ws.Range('I3').Value = 0.1827267341390226
$ws.Range('J3').Value = 0.1827267341390226
$ws.Range('K3').Value = 3
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = 4.349453666666666
$ws.Range('N3').Value = 13.048361
$ws.Range('O3').Value = 0.2460650667951621
$ws.Range('P3').Value = 0.2460650667951621
$ws.Range('Q3').Value = 151.0999986327317
$ws.Range('R3').Value = 1359.899987694585
$ws.Range('S3').Value = 0.04496266604118043
$ws.Range('T3').Value = 0.04496266604118043

# Row 4
$ws.Range('A4').Value = 'ECs'
$ws.Range('B4').Value = 'C3'
$ws.Range('C4').Value = 'Cd46'
$ws.Range('D4').Value = 'sCs'
$ws.Range('E4').Value = 2
$ws.Range('F4').Value = 0.6666666666666666
$ws.Range('G4').Value = 34.739995
$ws.Range('H4').Value = 104.219985
$ws.Range('I4').Value = 0.1827267341390226
$ws.Range('J4').Value = 0.1827267341390226
$ws.Range('K4').Value = 3
$ws.Range('L4').Value = 1
$ws.Range('M4').Value = 2.45022
$ws.Range('N4').Value = 7.35066
$ws.Range('O4').Value = 0.138618225222963
$ws.Range('P4').Value = 0.138618225222963
$ws.Range('Q4').Value = 85.1206305489
$ws.Range('R4').Value = 766.0856749400999
$ws.Range('S4').Value = 0.02532925558713952
$ws.Range('T4').Value = 0.02532925558713952

# Row 5
$ws.Range('A5').Value = 'FAPs'
$ws.Range('B5').Value = 'C3'
$ws.Range('C5').Value = 'Cd46'
$ws.Range('D5').Value = 'ECs'
$ws.Range('E5').Value = 3
$ws.Range('F5').Value = 1
$ws.Range('G5').Value = 154.8642143333334
$ws.Range('H5').Value = 464.5926430000001
$ws.Range('I5').Value = 0.8145606273154508
$ws.Range('J5').Value = 0.8145606273154508
$ws.Range('K5').Value = 3
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = 10.876357
$ws.Range('N5').Value = 32.629071
$ws.Range('O5').Value = 0.6153167079818751
$ws.Range('P5').Value = 0.615316707981875
$ws.Range('Q5').Value = 1684.35848161385
$ws.Range('R5').Value = 15159.22633452465
$ws.Range('S5').Value = 0.5012127636513942
$ws.Range('T5').Value = 0.5012127636513941

# Row 6
$ws.Range('A6').Value = 'FAPs'
$ws.Range('B6').Value = 'C3'
$ws.Range('C6').Value = 'Cd46'
$ws.Range('D6').Value = 'FAPs'
$ws.Range('E6').Value = 3
$ws.Range('F6').Value = 1
$ws.Range('G6').Value = 154.8642143333334
$ws.Range('H6').Value = 464.5926430000001
$ws.Range('I6').Value = 0.8145606273154508
$ws.Range('J6').Value = 0.8145606273154508
$ws.Range('K6').Value = 3
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = 4.349453666666666
$ws.Range('N6').Value = 13.048361
$ws.Range('O6').Value = 0.2460650667951621
$ws.Range('P6').Value = 0.2460650667951621
$ws.Range('Q6').Value = 673.5747248675693
$ws.Range('R6').Value = 6062.172523808124
$ws.Range('S6').Value = 0.2004349151690856
$ws.Range('T6').Value = 0.2004349151690855

# Row 7
$ws.Range('A7').Value = 'FAPs'
$ws.Range('B7').Value = 'C3'
$ws.Range('C7').Value = 'Cd46'
$ws.Range('D7').Value = 'sCs'
$ws.Range('E7').Value = 3
$ws.Range('F7').Value = 1
$ws.Range('G7').Value = 154.8642143333334
$ws.Range('H7').Value = 464.5926430000001
$ws.Range('I7').Value = 0.8145606273154508
$ws.Range('J7').Value = 0.8145606273154508
$ws.Range('K7').Value = 3
$ws.Range('L7').Value = 1
$ws.Range('M7').Value = 2.45022
$ws.Range('N7').Value = 7.35066
$ws.Range('O7').Value = 0.138618225222963
$ws.Range('P7').Value = 0.138618225222963
$ws.Range('Q7').Value = 379.4513952438201
$ws.Range('R7').Value = 3415.06255719438
$ws.Range('S7').Value = 0.1129129484949712
$ws.Range('T7').Value = 0.1129129484949711

# Row 8
$ws.Range('A8').Value = 'sCs'
$ws.Range('B8').Value = 'C3'
$ws.Range('C8').Value = 'Cd46'
$ws.Range('D8').Value = 'ECs'
$ws.Range('E8').Value = 3
$ws.Range('F8').Value = 1
$ws.Range('G8').Value = 0.5157266666666667
$ws.Range('H8').Value = 1.54718
$ws.Range('I8').Value = 0.002712638545526686
$ws.Range('J8').Value = 0.002712638545526686
$ws.Range('K8').Value = 3
$ws.Range('L8').Value = 1
$ws.Range('M8').Value = 10.876357
$ws.Range('N8').Value = 32.629071
$ws.Range('O8').Value = 0.6153167079818751
$ws.Range('P8').Value = 0.615316707981875
$ws.Range('Q8').Value = 5.609227341086666
$ws.Range('R8').Value = 50.48304606977999
$ws.Range('S8').Value = 0.001669131819778222
$ws.Range('T8').Value = 0.001669131819778222

# Row 9
$ws.Range('A9').Value = 'sCs'
$ws.Range('B9').Value = 'C3'
$ws.Range('C9').Value = 'Cd46'
$ws.Range('D9').Value = 'FAPs'
$ws.Range('E9').Value = 3
$ws.Range('F9').Value = 1
$ws.Range('G9').Value = 0.5157266666666667
$ws.Range('H9').Value = 1.54718
$ws.Range('I9').Value = 0.002712638545526686
$ws.Range('J9').Value = 0.002712638545526686
$ws.Range('K9').Value = 3
$ws.Range('L9').Value = 1
$ws.Range('M9').Value = 4.349453666666666
$ws.Range('N9').Value = 13.048361
$ws.Range('O9').Value = 0.2460650667951621
$ws.Range('P9').Value = 0.2460650667951621
$ws.Range('Q9').Value = 2.243129241331111
$ws.Range('R9').Value = 20.18816317198
$ws.Range('S9').Value = 0.0006674855848961554
$ws.Range('T9').Value = 0.0006674855848961554

# Row 10
$ws.Range('A10').Value = 'sCs'
$ws.Range('B10').Value = 'C3'
$ws.Range('C10').Value = 'Cd46'
$ws.Range('D10').Value = 'sCs'
$ws.Range('E10').Value = 3
$ws.Range('F10').Value = 1
$ws.Range('G10').Value = 0.5157266666666667
$ws.Range('H10').Value = 1.54718
$ws.Range('I10').Value = 0.002712638545526686
$ws.Range('J10').Value = 0.002712638545526686
$ws.Range('K10').Value = 3
$ws.Range('L10').Value = 1
$ws.Range('M10').Value = 2.45022
$ws.Range('N10').Value = 7.35066
$ws.Range('O10').Value = 0.138618225222963
$ws.Range('P10').Value = 0.138618225222963
$ws.Range('Q10').Value = 1.2636437932
$ws.Range('R10').Value = 11.3727941388
$ws.Range('S10').Value = 0.0003760211408523089
$ws.Range('T10').Value = 0.0003760211408523089
